# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ciudades")

# Swap the Almeria/Lugo rows (A47 <-> A48) along with their "Casos activos" (C column) values
$ws.Range("A47").Value = "Lugo"
$ws.Range("A48").Value = "Almeria"

$ws.Range("C47").Value = 5
$ws.Range("C48").Value = 72

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 22 de Marzo de 2020 a las 23:46"
